$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: LDPE / Aspergillus
$ws.Range("A8").Value = "LDPE"
$ws.Range("B8").Value = "Aspergillus"
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = 80
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 150
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 200
$ws.Range("I8").Value = 0.12
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = "Ahmed et al. 2025 (LDPE weight loss); fungal dye kinetic studies"

# Row 9: Phenol / Escherichia coli
$ws.Range("A9").Value = "Phenol"
$ws.Range("B9").Value = "Escherichia coli"
$ws.Range("C9").Value = 30
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 1.9
$ws.Range("F9").Value = 200
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 500
$ws.Range("I9").Value = 0.25
$ws.Range("J9").Value = 45
$ws.Range("K9").Value = "Sokół 1988; Haldane fits for phenol degradation by E. coli and Pseudomonas"

# Column B width (diff shows width 13.44140625 customWidth for col B)
$ws.Columns("B").ColumnWidth = 12.6

# Select K9 as active cell (matches sheetView selection)
$ws.Range("K9").Select()
